$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.945
$ws.Range("D4").Value = -7.523000000000001
$ws.Range("D6").Value = -7.906000000000001
$ws.Range("A9").Value = -20.783
$ws.Range("D10").Value = -7.592999999999999
$ws.Range("B11").Value = 7.215000000000001
$ws.Range("D11").Value = -8.294
$ws.Range("E12").Value = 12.92
$ws.Range("E17").Value = 13.604
$ws.Range("A18").Value = -21.825
$ws.Range("E19").Value = 12.759
$ws.Range("A20").Value = -21.738
$ws.Range("C21").Value = -12.031
$ws.Range("D21").Value = -7.717000000000001
